# InvoiceUpload.xlsx - add itemsDiscount / netTotal columns and T2(Tbl01) tax column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Make room for the two new "itemsDiscount" / "netTotal" columns between
#    the existing "salesTotal" (K) and "total" (old L) columns. The old
#    L/M/N ("total" / "T1(V009)" / "T1(V001)") shift to N/O/P, and the new
#    "T2(Tbl01)" tax column is written directly into the first free column Q.
# ---------------------------------------------------------------------------
$ws.Range("L1:M1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2) Row 1 - internal (English) header keys
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "itemsDiscount"
$ws.Range("M1").Value = "netTotal"
$ws.Range("Q1").Value = "T2(Tbl01)"

# ---------------------------------------------------------------------------
# 3) Row 2 - Arabic display headers
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "خصم الأصناف"
$ws.Range("M2").Value = "الصافى (بعد الخصم قبل الضريبة)"
$ws.Range("Q2").Value = "T2(Tbl01)"

# ---------------------------------------------------------------------------
# 4) Row 3 - sample data row
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 500
$ws.Range("C3").Value = 30
$ws.Range("E3").Value = 4690
$ws.Range("F3").Value = "EG-237791390-QTC1012"
$ws.Range("H3").Value = 11
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 1000
$ws.Range("N3").Value = 1145
$ws.Range("O3").Value = 140
$ws.Range("Q3").Value = 5

# ---------------------------------------------------------------------------
# 5) Column widths for the two newly inserted columns (best-fit autosize
#    performed by Excel on save).
# ---------------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 13.85546875
$ws.Columns.Item(13).ColumnWidth = 13.85546875

# ---------------------------------------------------------------------------
# 6) View - move the selection to the newly filled-in discount cell and
#    scroll the sheet so it is visible.
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$ws.Range("O3").Select()
